$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. Heading3 title change
Replace-Text "Mechanical Properties" "Physical Metallurgy"

# 2. Activation date
Replace-Text "Ativação: 01/01/2024" "Ativação: 01/01/2025"

# 3. Objetivos paragraph
$oldObjetivos = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências no desenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de materiais e a redução de ocorrência de falhas estruturais. Para tanto, a disciplina estabelece correlações com outras do curso de Engenharia de Materiais como LOM3013 – Ciência dos Materiais, LOM3057 – Introdução aos Materiais Poliméricos, LOM3032 - Cerâmica Física e LOM3011- Ensaios Mecânicos. Desta forma, são apresentadas a correlação entre propriedades e microestrutura de materiais para aplicações em Engenharia permitindo aos alunos a prática da redação científica e da busca bibliográfica para incentivar a solução de problemas em engenharia."
$newObjetivos = "Esta disciplina faz parte da formação do engenheiro de materiais e têm como objetivo gerar competências nodesenvolvimento de projetos seguros de equipamentos e componentes estruturais com o uso eficiente de metais e a redução de ocorrência de falhas estruturais baseado no trinômio propriedades, estrutura metalúrgica e processamento metalúrgico dos metais aplicado a engenharia permitindo aos alunos a prática da redação científica e da busca de projetos para incentivar a solução de problemas em engenharia."
Replace-Text $oldObjetivos $newObjetivos

# 4. Add new docente entry before Maria Ismenia
$found = $d.Content
$found.Find.Execute("7459752 - Maria Ismenia Sodero Toledo Faria") | Out-Null
$insertRange = $found.Duplicate
$insertRange.Collapse(1)
$insertRange.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert`v")

# 5. Programa resumido paragraph
$oldResumido = "1. Introdução ao conceito de propriedades mecânicas. 2. Elasticidade e Mecanismos de deformação plástica. 3. Teoria das discordâncias. 4.Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos. 6. Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos. 7. Influência da temperatura no comportamento mecânico de materiais. 8. Introdução básica à análise de falhas de materiais dúcteis e frágeis."
$newResumido = "1. Introdução ao conceito de propriedades mecânicas. 2. Elasticidade e mecanismos de deformação plástica. 3. Teoria das discordâncias. 4.Mecanismos de endurecimento. 5. Comportamento mecânico dos materiais metálicos. 6. Estudo comparativo de propriedades mecânicas de materiais metálicos. 7. Influência da temperatura no comportamento mecânico dos metais. 8. Introdução básica à análise de falhas de metais dúcteis e frágeis."
Replace-Text $oldResumido $newResumido

# 6. Programa paragraph
$oldPrograma = "1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais. Comportamento elástico e plástico de metais e ligas. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Tratamentos térmicos em aços. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência, fadiga de alto ciclo e propagação de trincas por fadiga. Impacto e a transição dúctil-frágil. 6. COMPORTAMENTO MECÂNICO DE MATERIAIS CERÂMICOS E POLIMÉRICOS: Estudo comparativo de propriedades mecânicas de materiais metálicos, cerâmicos e poliméricos 7. Influência da temperatura sobre o comportamento mecânico de materiais. Aspectos básicos  da  análise de falhas em materiais metálicos, cerâmicos e poliméricos."
$newPrograma = "Programa1.INTRODUÇÃO AO CONCEITO DE PROPRIEDADES MECÂNICAS: Conceitos e relações entre microestrutura e propriedades mecânicas de materiais metálicos. Comportamento elástico e plástico de metais suas ligas e materiais não ferrosos. 2. MECANISMOS DE DEFORMAÇÃO PLÁSTICA: Sistemas de deslizamento e movimentação de discordâncias. Deformação por maclação. Movimento relativo de grãos. Difusão. 3. TEORIA DAS DISCORDÂNCIAS: Classificação, observação e fontes de discordâncias. Multiplicação e interação de discordâncias. Forças entre discordâncias. Forças atuantes sobre discordâncias. Campos de tensão e energia. Energia de falha de empilhamento. Mecanismos de escalagem, deslizamento com desvio e empilhamento de discordâncias. Subestruturas de discordâncias. 4. MECANISMOS DE ENDURECIMENTO: Endurecimento por deformação plástica: Encruamento. Aumento da resistência devido aos contornos de grão. Relação de Hall-Petch. Endurecimento por solução sólida. Endurecimento por precipitação. Aços comuns e especiais. Estudo de ligas não metálicas. Tratamentos térmicos em aços e ligas especiais. 5. COMPORTAMENTO MECÂNICO DOS MATERIAIS METÁLICOS: Relação entre microestrutura e propriedades. Análise das propriedades em função de solicitações estáticas e cíclicas. Propriedades em tração uniaxial, fluência. Impacto e a transição dúctil-frágil. 6. Influência da temperatura sobre o comportamento mecânico dos metais. Aspectos básicos da análise de falhas em materiais metálicos."
Replace-Text $oldPrograma $newPrograma

# 7. Norma de recuperação text
$oldRecuperacao = "1. Meyers, M., Chawla, K. Mechanical Behavior of Materials. Ed. Cambridge University Press, 2009. 2. A. S. Lisbão, Estrutura e propriedades dos polímeros, EduFSCar, São Carlos, 2009. 3. T. H. Courtney, Mechanical Behavior of Materials, Waveland Press, 2005. 4. A. K. Bhargava, Engineering Materials: Polymers, Ceramics and Composites, PHI Learning Pvt. Ltd., 2012. 5.Dowling, E. M. Mechanical behavior of materials: engineering methods for deformation, fracture and fatigue. New Jersey, Prentice Hall, 2007. 6. Hull, D. Introduction to Dislocations, Pergamon Press, 1965. 7. Honeycombe, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1967. 8. Reed-Hill, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. 9. Van Vlack, L.H. Princípios de Ciência dos Materiais, Ed. Edgard Blucher Ltda., 1970. 10. Costa e Silva, A. L., Mei, P. R. Aços e Ligas especiais, Ed. Edgar Blücher, 2008. 11. Dieter, G.E. Metalurgia Mecânica, Ed. Guanabara Dois, 1986.  12. Callister, W. Ciência e engenharia dos materiais: Uma introdução, Rio de Janeiro, Livros Técnicos e Científicos, 2008. 13. Brooks, C. R., Choudhury, A. Metallurgical Failure Analysis, Ed. McGraw-Hill, 1993."
$newRecuperacao = "A recuperação será composta por uma única prova (PR) abrangendo toda a matéria ministrada ao longo do semestre. A Média final (MF) será computada pela relação:  MF = (NF + PR)/2."
Replace-Text $oldRecuperacao $newRecuperacao
